$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 4.9
$ws.Range("Q2").Value = 1.68
$ws.Range("S2").Value = 2.66
$ws.Range("AH2").Value = 19.5
$ws.Range("G3").Value = 2.64
$ws.Range("Q3").Value = 2.56
$ws.Range("G4").Value = 2.54
$ws.Range("I4").Value = 3.2
$ws.Range("N4").Value = 5.5
$ws.Range("R4").Value = 1.63
$ws.Range("S4").Value = 2.3
$ws.Range("T4").Value = 1.5
$ws.Range("U4").Value = 2.58
$ws.Range("V4").Value = 1.46
$ws.Range("Y4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 1000
$ws.Range("AN4").Value = 12.5
$ws.Range("AO4").Value = 1000
$ws.Range("G6").Value = 1.7
$ws.Range("M6").Value = 1.04
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.56
$ws.Range("T6").Value = 1.81
$ws.Range("U6").Value = 1.96
$ws.Range("W6").Value = 2.42
$ws.Range("H7").Value = 2.94
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.22
$ws.Range("S7").Value = 2.56
$ws.Range("H8").Value = 7.2
$ws.Range("K8").Value = 6
$ws.Range("H9").Value = 7
$ws.Range("I9").Value = 7.4
$ws.Range("N9").Value = 2.92
$ws.Range("T9").Value = 2.38
$ws.Range("W9").Value = 2.44
$ws.Range("H10").Value = 24
$ws.Range("P10").Value = 2.98
$ws.Range("Q10").Value = 1.38
$ws.Range("R10").Value = 1.78
$ws.Range("S10").Value = 1.98
$ws.Range("T10").Value = 2.7
$ws.Range("AL10").Value = 80
